{"js": "// Office.js (Word JavaScript API) edit script.\n// 1) Merge the two runs in \"Swing trading (1+ minute trades\" + \")\" into one run.\n// 2) Merge the two runs in \"Trend trading (buying / selling the trend\" + \")\" into one run.\n// 3) After the paragraph ending \"... this algorithm will potentially be \", insert\n//    two empty paragraphs followed by a paragraph containing the new question text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1 & 2: merge the split \"(...)\" runs back into a single run each. ---\n// Re-writing the paragraph's text via insertText(..., \"Replace\") collapses all\n// of the paragraph's runs into a single run carrying the full text, which is\n// exactly the shape the diff shows (two <w:r> -> one <w:r>).\nconst swingPara = paragraphs.items.filter(p => p.text.indexOf(\"Swing trading (1+ minute trades\") !== -1)[0];\nif (swingPara) {\n  swingPara.insertText(\"Swing trading (1+ minute trades)\", \"Replace\");\n}\n\nconst trendPara = paragraphs.items.filter(p => p.text.indexOf(\"Trend trading (buying / selling the trend\") !== -1)[0];\nif (trendPara) {\n  trendPara.insertText(\"Trend trading (buying / selling the trend)\", \"Replace\");\n}\n\n// --- 3: insert the two blank paragraphs + the new question paragraph. ---\nconst anchorPara = paragraphs.items.filter(p => p.text.indexOf(\"this algorithm will potentially be\") !== -1)[0];\nif (anchorPara) {\n  const blank1 = anchorPara.insertParagraph(undefined, \"After\");\n  const blank2 = blank1.insertParagraph(undefined, \"After\");\n  blank2.insertParagraph(\n    \"Should the initial plan be for an AI that trades based off charts and then introduce the learning capabilities first?\",\n    \"After\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# 1) Merge the two runs in \"Swing trading (1+ minute trades\" + \")\" into one run.\n# 2) Merge the two runs in \"Trend trading (buying / selling the trend\" + \")\" into one run.\n# 3) After the paragraph ending \"... this algorithm will potentially be \", insert\n#    two empty paragraphs followed by a paragraph containing the new question text.\n\n$d = $word.ActiveDocument\n\n# --- helper: locate the paragraph containing a unique snippet of text ---\nfunction Get-ParagraphContaining($doc, [string]$snippet) {\n    $find = $doc.Content\n    $find.Find.Execute($snippet) | Out-Null\n    return $find.Paragraphs(1)\n}\n\n# --- 1 & 2: merge the split \"(...)\" runs back into a single run each. ---\n# Re-writing the paragraph's own Range.Text collapses every run it owns into a\n# single run carrying the full text, matching the two-runs -> one-run shape in\n# the diff. A freshly-fetched Range (via $d.Range(start,end)) is used instead\n# of the paragraph's live Range object so the write isn't clipped to the first\n# run.\n$swingPara = Get-ParagraphContaining $d \"Swing trading (1+ minute trades\"\n$swingRange = $d.Range($swingPara.Range.Start, $swingPara.Range.End)\n$swingRange.Text = \"Swing trading (1+ minute trades)\"\n\n$trendPara = Get-ParagraphContaining $d \"Trend trading (buying / selling the trend\"\n$trendRange = $d.Range($trendPara.Range.Start, $trendPara.Range.End)\n$trendRange.Text = \"Trend trading (buying / selling the trend)\"\n\n# --- 3: insert the two blank paragraphs + the new question paragraph. ---\n$anchorPara = Get-ParagraphContaining $d \"this algorithm will potentially be\"\n\n$anchorPara.Range.InsertParagraphAfter()\n$blank1Start = $anchorPara.Range.End\n$blank1 = $d.Range($blank1Start, $blank1Start).Paragraphs(1)\n\n$blank1.Range.InsertParagraphAfter()\n$blank2Start = $blank1.Range.End\n$blank2 = $d.Range($blank2Start, $blank2Start).Paragraphs(1)\n\n$blank2.Range.InsertParagraphAfter()\n$textStart = $blank2.Range.End\n$questionPara = $d.Range($textStart, $textStart).Paragraphs(1)\n$questionPara.Range.Text = \"Should the initial plan be for an AI that trades based off charts and then introduce the learning capabilities first?\"\n"}
